$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from H1 into the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF)
$data = @(
    @(8, 9),
    @(8, 8),
    @(2, 3),
    @(9, 9),
    @(5, 6),
    @(4, 5),
    @(8, 9),
    @(8, 9),
    @(2, 3),
    @(8, 8),
    @(8, 8),
    @(1, 1),
    @(7, 7),
    @(7, 8),
    @(1, 1),
    @(9, 9),
    @(4, 5),
    @(5, 5),
    @(7, 7),
    @(2, 2),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
